$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(6).Delete()
$ws.Range("F2").Copy($ws.Range("F42"))
$ws.Range("F2").Copy($ws.Range("F43"))
$ws.Range("F42").Value = 0.15566
$ws.Range("F43").Value = 0.15566
